$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Summary" (sheet1): update aggregate stats now that trade #9 closed.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200        # Current Capital
$summary.Range("B4").Value = 0           # Total P&L $
$summary.Range("B5").Value = 0           # Total P&L %
$summary.Range("B6").Value = 9           # Total Trades
$summary.Range("B8").Value = 3           # Losing Trades
$summary.Range("B9").Value = 44.44       # Win Rate %

# ---------------------------------------------------------------------------
# Sheet "Strategy Status" (sheet2): update the MarketMaking strategy row (4).
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100          # Capital
$status.Range("D4").Value = 9            # Trades
$status.Range("E4").Value = 0            # P&L $
$status.Range("F4").Value = 0            # P&L %
$status.Range("G4").Value = 44.44        # Win Rate %

# ---------------------------------------------------------------------------
# New trade #9 row, appended to both "All Trades" and "MarketMaking" sheets.
# ---------------------------------------------------------------------------
$newRow = @(9, "2026-02-17", "08:08:20", "MarketMaking", "UP", 0.42, 0.33, "CLOSED", -21.4286, -0.09, 100, 0, 0, 0.6, "Normal spread capture: 19600 bps", "early_exit", 0.13)

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $row = 10
    $col = 1
    foreach ($value in $newRow) {
        $cell = $ws.Cells.Item($row, $col)
        if ($col -eq 2) {
            # The "Date" column holds plain text like "2026-02-17" in the
            # source data (not an actual date cell). Force text formatting
            # before assignment so Excel doesn't auto-convert it into a
            # date serial number, then drop the temporary format again so
            # no stray style sticks around on the cell.
            $cell.NumberFormat = "@"
            $cell.Value = $value
            $cell.ClearFormats()
        } else {
            $cell.Value = $value
        }
        $col = $col + 1
    }
}
